$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 613.1111
$ws.Range("I20").Value = 684.7143
$ws.Range("J20").Value = 362.5
$ws.Range("K20").Value = 684.7143
$ws.Range("L20").Value = 362.5
$ws.Range("M20").Value = -454.7143
$ws.Range("N20").Value = -822.5

# Row 35
$ws.Range("H35").Value = 613.1111
$ws.Range("I35").Value = 684.7143
$ws.Range("J35").Value = 362.5
$ws.Range("K35").Value = 684.7143
$ws.Range("L35").Value = 362.5
$ws.Range("M35").Value = -305.7143
$ws.Range("N35").Value = -1120.5

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 80
$ws.Range("H80").Value = 601.7143
$ws.Range("I80").Value = 601.7143
$ws.Range("K80").Value = 1805.1429
$ws.Range("M80").Value = -807.1428999999998

# Row 83
$ws.Range("H83").Value = 601.7143
$ws.Range("I83").Value = 601.7143
$ws.Range("K83").Value = 5415.428699999999
$ws.Range("M83").Value = -423.4286999999995

# Row 133
$ws.Range("H133").Value = 22999
$ws.Range("J133").Value = 22999
$ws.Range("L133").Value = 22999
$ws.Range("N133").Value = -33119

# Row 138
$ws.Range("H138").Value = 1294.7916
$ws.Range("I138").Value = 462.27274
$ws.Range("K138").Value = 1386.81822
$ws.Range("M138").Value = 3753.18178


$ws = $wb.Worksheets.Item("ARM")
# Row 80
$ws.Range("H80").Value = 39000
$ws.Range("J80").Value = 39000
$ws.Range("L80").Value = 39000
$ws.Range("N80").Value = -40996

# Row 83
$ws.Range("H83").Value = 39000
$ws.Range("J83").Value = 39000
$ws.Range("L83").Value = 117000
$ws.Range("N83").Value = -126984

# Row 132
$ws.Range("H132").Value = 1898.75
$ws.Range("I132").Value = 1898.75
$ws.Range("K132").Value = 5696.25
$ws.Range("M132").Value = -3166.25


$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 11500
$ws.Range("J35").Value = 11500
$ws.Range("L35").Value = 11500
$ws.Range("N35").Value = -12120

# Row 48
$ws.Range("H48").Value = 199999
$ws.Range("J48").Value = 199999
$ws.Range("L48").Value = 199999
$ws.Range("N48").Value = -200829

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 82
$ws.Range("H82").Value = 31494.625
$ws.Range("I82").Value = 5978.5
$ws.Range("K82").Value = 5978.5
$ws.Range("M82").Value = -5595.5

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 85
$ws.Range("H85").Value = 31494.625
$ws.Range("I85").Value = 5978.5
$ws.Range("K85").Value = 5978.5
$ws.Range("M85").Value = -4652.5

# Row 107
$ws.Range("H107").Value = 903
$ws.Range("I107").Value = 973.25
$ws.Range("K107").Value = 973.25
$ws.Range("M107").Value = 946.75


$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 19983
$ws.Range("I41").Value = 4898
$ws.Range("J41").Value = 23000
$ws.Range("K41").Value = 4898
$ws.Range("L41").Value = 23000
$ws.Range("M41").Value = -4470
$ws.Range("N41").Value = -23856

# Row 50
$ws.Range("H50").Value = 29996.875
$ws.Range("J50").Value = 29996.875
$ws.Range("L50").Value = 29996.875
$ws.Range("N50").Value = -31246.875

# Row 59
$ws.Range("H59").Value = 34997.9
$ws.Range("J59").Value = 34997.9
$ws.Range("L59").Value = 34997.9
$ws.Range("N59").Value = -37287.9

# Row 60
$ws.Range("H60").Value = 20737.8
$ws.Range("I60").Value = 10797.667
$ws.Range("J60").Value = 24997.857
$ws.Range("K60").Value = 10797.667
$ws.Range("L60").Value = 24997.857
$ws.Range("M60").Value = -10286.667
$ws.Range("N60").Value = -26019.857

# Row 68
$ws.Range("H68").Value = 39010.4
$ws.Range("I68").Value = 20268
$ws.Range("J68").Value = 39996.844
$ws.Range("K68").Value = 20268
$ws.Range("L68").Value = 39996.844
$ws.Range("M68").Value = -19519
$ws.Range("N68").Value = -41494.844

# Row 71
$ws.Range("H71").Value = 39010.4
$ws.Range("I71").Value = 20268
$ws.Range("J71").Value = 39996.844
$ws.Range("K71").Value = 60804
$ws.Range("L71").Value = 119990.532
$ws.Range("M71").Value = -57060
$ws.Range("N71").Value = -127478.532

# Row 74
$ws.Range("H74").Value = 39747.07
$ws.Range("I74").Value = 36494
$ws.Range("J74").Value = 39997.31
$ws.Range("K74").Value = 36494
$ws.Range("L74").Value = 39997.31
$ws.Range("M74").Value = -35620
$ws.Range("N74").Value = -41745.31

# Row 77
$ws.Range("H77").Value = 39747.07
$ws.Range("I77").Value = 36494
$ws.Range("J77").Value = 39997.31
$ws.Range("K77").Value = 109482
$ws.Range("L77").Value = 119991.93
$ws.Range("M77").Value = -105114
$ws.Range("N77").Value = -128727.93

# Row 107
$ws.Range("H107").Value = 1086.2222
$ws.Range("I107").Value = 1129.25
$ws.Range("K107").Value = 1129.25
$ws.Range("M107").Value = 790.75

# Row 134
$ws.Range("H134").Value = 4909.0835
$ws.Range("I134").Value = 4512
$ws.Range("J134").Value = 5703.25
$ws.Range("K134").Value = 13536
$ws.Range("L134").Value = 17109.75
$ws.Range("M134").Value = -11001
$ws.Range("N134").Value = -22179.75


$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 13934047
$ws.Range("J122").Value = 6000
$ws.Range("L122").Value = 18000
$ws.Range("N122").Value = -22900

# Row 132
$ws.Range("H132").Value = 2396.875
$ws.Range("J132").Value = 2674.25
$ws.Range("L132").Value = 8022.75
$ws.Range("N132").Value = -13082.75


$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 768.1667
$ws.Range("I81").Value = 721.125
$ws.Range("J81").Value = 862.25
$ws.Range("K81").Value = 1442.25
$ws.Range("L81").Value = 1724.5
$ws.Range("M81").Value = -381.25
$ws.Range("N81").Value = -3846.5

# Row 84
$ws.Range("H84").Value = 768.1667
$ws.Range("I84").Value = 721.125
$ws.Range("J84").Value = 862.25
$ws.Range("K84").Value = 7211.25
$ws.Range("L84").Value = 8622.5
$ws.Range("M84").Value = -1907.25
$ws.Range("N84").Value = -19230.5

# Row 107
$ws.Range("H107").Value = 1490
$ws.Range("I107").Value = 1490
$ws.Range("K107").Value = 4470
$ws.Range("M107").Value = -2550

# Row 126
$ws.Range("H126").Value = 2292.3333
$ws.Range("I126").Value = 1990.2727
$ws.Range("J126").Value = 3123
$ws.Range("K126").Value = 5970.8181
$ws.Range("L126").Value = 9369
$ws.Range("M126").Value = -3500.8181
$ws.Range("N126").Value = -14309

# Row 137
$ws.Range("H137").Value = 49999
$ws.Range("J137").Value = 49999
$ws.Range("L137").Value = 49999
$ws.Range("N137").Value = -60199

